# "safeareaview and status bar"
#
# Appends a new slide (slide 11 / sldId 266) at the end of the deck using the
# same "Title Only" layout family already used by the other content slides,
# and gives it the title "5. React Native Components" — continuing the
# deck's numbered section pattern (1, 1.1, 1.2, 1.3, 2, 3, 3.2, 4, 5).

$p = $ppt.ActivePresentation

# Append a new slide after the current last slide.
# ppLayoutTitleOnly (11) -> single Title placeholder, no Content placeholder,
# matching the target slide which only contains a Title shape.
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 11)

# Set the slide's title text.
$slide.Shapes.Item(1).TextFrame.TextRange.Text = "5. React Native Components"
